$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.976.27'
$ws.Range("E2").Value = '  +4.79%  '
$ws.Range("D3").Value = '3.080.44'
$ws.Range("E3").Value = '  +2.88%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.56'
$ws.Range("E5").Value = '  +3.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.40'
$ws.Range("E6").Value = '  +2.13%  '
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D8").Value = '3.072.34'
$ws.Range("E8").Value = '  +3.16%  '
$ws.Range("E9").Value = '  +1.72%  '
$ws.Range("E10").Value = '  +5.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.58'
$ws.Range("E11").Value = '  +8.85%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.468'
$ws.Range("E12").Value = '  +2.87%  '
$ws.Range("E13").Value = '  +4.41%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.41'
$ws.Range("E14").Value = '  +5.11%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.124'
$ws.Range("E15").Value = '  +0.32%  '
$ws.Range("D16").Value = '3.590.31'
$ws.Range("E16").Value = '  +2.91%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.31'
$ws.Range("E17").Value = '  +3.20%  '
$ws.Range("D18").Value = '3.077.97'
$ws.Range("E18").Value = '  +2.85%  '
$ws.Range("D19").Value = '61.881.01'
$ws.Range("E19").Value = '  +4.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '448.39'
$ws.Range("E20").Value = '  +4.74%  '
$ws.Range("E21").Value = '  +2.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.731'
$ws.Range("E22").Value = '  +2.37%  '
$ws.Range("E23").Value = '  +5.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.85'
$ws.Range("E24").Value = '  +3.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.05'
$ws.Range("E25").Value = '  +1.60%  '
$ws.Range("E26").Value = '  +0.11%  '
$ws.Range("E27").Value = '  +6.31%  '
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("E29").Value = '  +4.98%  '
$ws.Range("E30").Value = '  +5.55%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.78'
$ws.Range("E31").Value = '  +10.84%  '
$ws.Range("E32").Value = '  +14.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.81'
$ws.Range("E33").Value = '  +4.18%  '
$ws.Range("E34").Value = '  +4.98%  '
$ws.Range("D35").Value = '0.0₃0797'
$ws.Range("E35").Value = '  +2.94%  '
$ws.Range("E36").Value = '  +3.27%  '
$ws.Range("E37").Value = '  +5.64%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '50.36'
$ws.Range("E38").Value = '  +2.14%  '
$ws.Range("E39").Value = '  +10.21%  '
$ws.Range("E40").Value = '  +2.42%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '429.77'
$ws.Range("E41").Value = '  +6.99%  '
$ws.Range("E42").Value = '  +5.73%  '
$ws.Range("D43").Value = '2.870.68'
$ws.Range("E43").Value = '  +3.36%  '
$ws.Range("E44").Value = '  +7.77%  '
$ws.Range("E45").Value = '  +1.29%  '
$ws.Range("E46").Value = '  +6.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '35.19'
$ws.Range("E47").Value = '  +4.08%  '
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.68'
$ws.Range("E49").Value = '  +1.75%  '
$ws.Range("E50").Value = '  +1.39%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.08'
$ws.Range("E51").Value = '  +2.70%  '
